# Re-upload of Sample_file.xlsx: refresh student data, drop the
# leftover "sample" sheet, and tidy up each sheet's contents.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ------------------------------------------------------------------
# 1) Remove the throwaway "sample" worksheet entirely.
# ------------------------------------------------------------------
$wsSample = $wb.Worksheets.Item("sample")
$wsSample.Delete()

# ------------------------------------------------------------------
# 2) "Data" sheet: add a proper header row above the existing record
#    and refresh a couple of the record's values.
# ------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")

$wsData.Rows.Item(1).Insert()

$wsData.Range("A1").Value = "Student No."
$wsData.Range("B1").Value = "Fullname"
$wsData.Range("C1").Value = "Email"
$wsData.Range("D1").Value = "Gender"
$wsData.Range("E1").Value = "Course"
$wsData.Range("F1").Value = "Contact No."
$wsData.Range("G1").Value = "Address"

$wsData.Range("B2").Value = "Jiro Miko S. Viñas"
$wsData.Range("G2").Value = "Calmar Homes Subdivision, Lucena City"

# Widen the columns that now hold longer text so everything is
# readable (best-fit-style sizing).
$wsData.Columns.Item(3).ColumnWidth = 24.5
$wsData.Columns.Item(6).ColumnWidth = 12
$wsData.Columns.Item(7).ColumnWidth = 36.333333333333336

# ------------------------------------------------------------------
# 3) "Old Students" sheet: refresh the sample login rows.
# ------------------------------------------------------------------
$wsOld = $wb.Worksheets.Item("Old Students")

$wsOld.Range("A2").Value = "022A-9661"
$wsOld.Range("B2").Value = "jiromiko"
$wsOld.Range("C2").Value = "BSIT-1A"

$wsOld.Range("A3").Value = "0231-1232"
$wsOld.Range("B3").Value = "aaaaa000"
$wsOld.Range("C3").Value = "1A"

$wsOld.Range("A4").Value = "0123-1234"
$wsOld.Range("B4").Value = "aaaaa000"
$wsOld.Range("C4").Value = "a"

# ------------------------------------------------------------------
# 4) "New Students" sheet: drop the placeholder data row, leaving
#    only the header behind.
# ------------------------------------------------------------------
$wsNew = $wb.Worksheets.Item("New Students")
$wsNew.Rows.Item(2).Delete()

# ------------------------------------------------------------------
# 5) Restore per-sheet selections, finishing on "New Students" so it
#    ends up the active tab.
# ------------------------------------------------------------------
$wsData.Range("A3:G3").Select()
$wsOld.Range("A3:C4").Select()
$wsNew.Activate()
$wsNew.Range("L12").Select()

Write-Output "done"
